$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.908.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.45%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.641.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.04%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5049'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '

$ws.Range('E7').Value = '  -0.45%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2567'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.30%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06390'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.59%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.82%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07802'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.06%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.658.77'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.25%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.278'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5428'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.47%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅7863'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.47%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.79'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.90%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.957.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.004'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.48%  '

$ws.Range('E19').Value = '  -2.71%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.394'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.31%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.954'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.42%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.986'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.29%  '

$ws.Range('E23').Value = '  -0.43%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.868'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.03%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1143'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.848'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.68%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.04%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.242'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04935'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.44%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.269'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.196'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.532'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.69%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.367'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.13%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8937'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.17%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.609'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.17%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.139.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.86%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5547'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.32%  '

$ws.Range('E39').Value = '  -0.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.004'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.671'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.16%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8189'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.79%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.54%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₈120'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.11%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.778.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4530'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.65%  '

$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05053'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.25%  '

$ws.Range('E50').Value = '  -0.05%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09507'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.03%  '
